$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column S by copying formatting from column R for the rows that need it ---

# Row 3: S3 should take on style of R3 (empty, bottom-border style)
$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)

# Row 4: S4 = 2022, same style as R4 (year header style)
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Row 5: update P5/Q5/R5 values, add S5 with style copied from R5
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("P5").Value = 23.111083656771282
$ws.Range("Q5").Value = 24.08077930418019
$ws.Range("R5").Value = 19.336931533747723
$ws.Range("S5").Value = 13.600365850576139

# Row 6: update P6/Q6/R6 values, add S6 with style copied from R6
$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("P6").Value = 14.322631450320875
$ws.Range("Q6").Value = 13.073459110725862
$ws.Range("R6").Value = 10.464141365743002
$ws.Range("S6").Value = 9.2742414863791556

# Row 7: update P7 value, add S7 with style copied from R7
$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("P7").Value = 23.612622725489956
$ws.Range("S7").Value = 17.303523954725925

# Row 8: add S8 with style copied from R8
$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$ws.Range("S8").Value = 205.5

# Clear clipboard marching-ants / copy mode
$excel.CutCopyMode = 0

# Update the saved selection to match the committed workbook state
$null = $ws.Range("Q15").Select()
